$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.772.91"
$ws.Range("E2").Value = "  -1.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.496.18"
$ws.Range("E3").Value = "  -3.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'579.45"
$ws.Range("E5").Value = "  -4.50%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'192.51"
$ws.Range("E6").Value = "  -3.58%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -2.44%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.485.09"
$ws.Range("E8").Value = "  -3.73%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.14%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -6.70%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -4.41%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'51.43"
$ws.Range("E12").Value = "  -4.49%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -6.58%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'9.15"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.051.97"
$ws.Range("E15").Value = "  -3.64%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "'647.05"
$ws.Range("E16").Value = "  -4.83%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.688.23"
$ws.Range("E17").Value = "  -1.79%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.495.10"
$ws.Range("E18").Value = "  -3.20%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'12.39"
$ws.Range("E19").Value = "  -4.58%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -1.90%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'18.31"
$ws.Range("E21").Value = "  -3.87%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.948"
$ws.Range("E22").Value = "  -5.19%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'18.02"
$ws.Range("E23").Value = "  -3.41%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'5.33"
$ws.Range("E24").Value = "  -1.18%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'98.76"
$ws.Range("E25").Value = "  -6.74%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "'4.29"
$ws.Range("E26").Value = "  -7.45%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  -4.64%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  -4.60%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'9.37"
$ws.Range("E29").Value = "  -5.19%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'32.67"
$ws.Range("E30").Value = "  -5.02%  "

# Row 31 - dogwifhat
$ws.Range("D31").Value = "'4.27"
$ws.Range("E31").Value = "  -7.70%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  -6.72%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "'11.66"
$ws.Range("E33").Value = "  -4.54%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -5.23%  "

# Row 35 - OKB
$ws.Range("D35").Value = "'61.21"
$ws.Range("E35").Value = "  -3.33%  "

# Row 36 - Bittensor
$ws.Range("D36").Value = "'563.60"
$ws.Range("E36").Value = "  +9.86%  "

# Row 37 - Maker
$ws.Range("D37").Value = "3.710.96"
$ws.Range("E37").Value = "  -6.02%  "

# Row 38 - Dai
$ws.Range("E38").Value = "  +0.22%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0790"
$ws.Range("E39").Value = "  -9.38%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "'3.61"
$ws.Range("E40").Value = "  +1.81%  "

# Row 41 - now CoreDAO (was Fetch.AI)
$ws.Range("B41").Value = "CoreDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D41").Value = "'3.75"
$ws.Range("E41").Value = "  +46.75%  "

# Row 42 - now Fetch.AI (was TheGraph)
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.91"
$ws.Range("E42").Value = "  -3.76%  "

# Row 43 - now TheGraph (was CoreDAO)
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.374"
$ws.Range("E43").Value = "  -3.81%  "

# Row 44 - Kaspa
$ws.Range("D44").Value = "'0.134"
$ws.Range("E44").Value = "  -2.23%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'34.28"
$ws.Range("E45").Value = "  -6.67%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "'0.0443"
$ws.Range("E46").Value = "  -3.73%  "

# Row 47 - ApeXProtocol
$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  -4.06%  "

# Row 48 - ThetaToken
$ws.Range("E48").Value = "  -7.93%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -4.48%  "

# Row 50 - FirstDigitalUSD
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.40%  "

# Row 51 - THORChain
$ws.Range("E51").Value = "  -5.40%  "
